$wb = $excel.ActiveWorkbook

# Sheets "展览" and "全部类型" both carry the same 31-row event listing.
# The refreshed scrape replaces the body (rows 2-26) with updated event
# data and drops the now-stale trailing rows (old rows 27-31), shrinking
# the used range from A1:I31 down to A1:I26.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Column B holds plain-text dates ("2024.02.14"); force text format
    # first so Excel does not auto-convert the strings into date serials.
    $ws.Range("B2:B26").NumberFormat = "@"

    $ws.Cells.Item(2, 2).Value = "2024.02.14"
    $ws.Cells.Item(2, 3).Value = "南昌·龙年动漫展"
    $ws.Cells.Item(2, 4).Value = "南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆"
    $ws.Cells.Item(2, 5).Value = "2024.02.14 10:00-02.15 18:00"
    $ws.Cells.Item(2, 6).Value = 1399
    $ws.Cells.Item(2, 7).Value = "已停售"
    $ws.Cells.Item(2, 8).Value = "https://show.bilibili.com/platform/detail.html?id=80525"
    $ws.Cells.Item(2, 9).Value = "//i2.hdslb.com/bfs/openplatform/202401/ezt7koZo1704444854691.jpeg"

    $ws.Cells.Item(3, 2).Value = "2024.02.15"
    $ws.Cells.Item(3, 3).Value = "赣州·明日方舟ONLY大炎新岁同好交流茶话会"
    $ws.Cells.Item(3, 4).Value = "南门口地一大道下沉广场 漫库书店"
    $ws.Cells.Item(3, 5).Value = "2024.02.15 11:00-02.15 18:00"
    $ws.Cells.Item(3, 6).Value = 169
    $ws.Cells.Item(3, 7).Value = "已停售"
    $ws.Cells.Item(3, 8).Value = "https://show.bilibili.com/platform/detail.html?id=78689"
    $ws.Cells.Item(3, 9).Value = "//i1.hdslb.com/bfs/openplatform/202311/T1Y8Iju31700621742031.png"

    $ws.Cells.Item(4, 2).Value = "2024.02.16"
    $ws.Cells.Item(4, 3).Value = "上高·星语动漫嘉年华"
    $ws.Cells.Item(4, 4).Value = "镜山大道2号 迎宾馆大酒店"
    $ws.Cells.Item(4, 5).Value = "2024.02.16 09:30-02.16 17:00"
    $ws.Cells.Item(4, 6).Value = 135
    $ws.Cells.Item(4, 7).Value = 40
    $ws.Cells.Item(4, 8).Value = "https://show.bilibili.com/platform/detail.html?id=80844"
    $ws.Cells.Item(4, 9).Value = "//i1.hdslb.com/bfs/openplatform/202401/QCJN9j8h1705306410081.png"

    $ws.Cells.Item(5, 2).Value = "2024.02.16"
    $ws.Cells.Item(5, 3).Value = "南昌·运动番only"
    $ws.Cells.Item(5, 4).Value = "南龙蟠街666号 融创茂"
    $ws.Cells.Item(5, 5).Value = "2024.02.16 10:00-02.16 17:00"
    $ws.Cells.Item(5, 6).Value = 301
    $ws.Cells.Item(5, 7).Value = 60
    $ws.Cells.Item(5, 8).Value = "https://show.bilibili.com/platform/detail.html?id=80757"
    $ws.Cells.Item(5, 9).Value = "//i2.hdslb.com/bfs/openplatform/202401/QXLfgq7f1706180123892.jpeg"

    $ws.Cells.Item(6, 2).Value = "2024.02.17"
    $ws.Cells.Item(6, 3).Value = "九江·ACD动漫游戏嘉年华02"
    $ws.Cells.Item(6, 4).Value = "九瑞大道与重庆路交汇处西南角 九江国际会展中心"
    $ws.Cells.Item(6, 5).Value = "2024.02.17 10:00-02.17 17:00"
    $ws.Cells.Item(6, 6).Value = 349
    $ws.Cells.Item(6, 7).Value = 55
    $ws.Cells.Item(6, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81055"
    $ws.Cells.Item(6, 9).Value = "//i0.hdslb.com/bfs/openplatform/202401/7BLpSOEZ1705574359625.jpeg"

    $ws.Cells.Item(7, 2).Value = "2024.02.17"
    $ws.Cells.Item(7, 3).Value = "江西·樟树静卿国风动漫文化展览会"
    $ws.Cells.Item(7, 4).Value = "樟树市杏佛路89号 银河国际酒店"
    $ws.Cells.Item(7, 5).Value = "2024.02.17 09:00-02.17 17:00"
    $ws.Cells.Item(7, 6).Value = 343
    $ws.Cells.Item(7, 7).Value = 40
    $ws.Cells.Item(7, 8).Value = "https://show.bilibili.com/platform/detail.html?id=80795"
    $ws.Cells.Item(7, 9).Value = "//i2.hdslb.com/bfs/openplatform/202401/DWQnrbtu1705044465383.jpeg"

    $ws.Cells.Item(8, 2).Value = "2024.02.17"
    $ws.Cells.Item(8, 3).Value = "赣州·第一届喵喵鱼动漫游戏展"
    $ws.Cells.Item(8, 4).Value = "105国道东100米 毅德国际会展中心"
    $ws.Cells.Item(8, 5).Value = "2024.02.17 09:30-02.18 16:00"
    $ws.Cells.Item(8, 6).Value = 1836
    $ws.Cells.Item(8, 7).Value = 50
    $ws.Cells.Item(8, 8).Value = "https://show.bilibili.com/platform/detail.html?id=78362"
    $ws.Cells.Item(8, 9).Value = "//i0.hdslb.com/bfs/openplatform/202311/KXRHxTLL1699521247861.png"

    $ws.Cells.Item(9, 2).Value = "2024.02.18"
    $ws.Cells.Item(9, 3).Value = "万载·第七届馨缘动漫文化展"
    $ws.Cells.Item(9, 4).Value = "康乐街道阳乐大道217号 龙凤大酒店"
    $ws.Cells.Item(9, 5).Value = "2024.02.18 09:30-02.18 17:00"
    $ws.Cells.Item(9, 6).Value = 75
    $ws.Cells.Item(9, 7).Value = 40
    $ws.Cells.Item(9, 8).Value = "https://show.bilibili.com/platform/detail.html?id=80971"
    $ws.Cells.Item(9, 9).Value = "//i1.hdslb.com/bfs/openplatform/202401/6ZDl6Oou1705487204077.png"

    $ws.Cells.Item(10, 2).Value = "2024.02.18"
    $ws.Cells.Item(10, 3).Value = "奉新·COP动漫游戏嘉年华1.0"
    $ws.Cells.Item(10, 4).Value = "应星北大道482号 金勺宴大酒店"
    $ws.Cells.Item(10, 5).Value = "2024.02.18 09:00-02.18 17:00"
    $ws.Cells.Item(10, 6).Value = 113
    $ws.Cells.Item(10, 7).Value = 30
    $ws.Cells.Item(10, 8).Value = "https://show.bilibili.com/platform/detail.html?id=78259"
    $ws.Cells.Item(10, 9).Value = "//i0.hdslb.com/bfs/openplatform/202311/yqw3kAkh1699597195072.jpeg"

    $ws.Cells.Item(11, 2).Value = "2024.02.20"
    $ws.Cells.Item(11, 3).Value = "江西·高安首届静卿国风动漫文化展览会"
    $ws.Cells.Item(11, 4).Value = "华林中路606号 华鼎国际大酒店"
    $ws.Cells.Item(11, 5).Value = "2024.02.20 09:00-02.20 17:00"
    $ws.Cells.Item(11, 6).Value = 187
    $ws.Cells.Item(11, 7).Value = 40
    $ws.Cells.Item(11, 8).Value = "https://show.bilibili.com/platform/detail.html?id=80785"
    $ws.Cells.Item(11, 9).Value = "//i0.hdslb.com/bfs/openplatform/202401/kcU6CEz91705040408216.jpeg"

    $ws.Cells.Item(12, 2).Value = "2024.02.23"
    $ws.Cells.Item(12, 3).Value = "上饶·囧喵喵次元国风动漫游戏展"
    $ws.Cells.Item(12, 4).Value = "春江北大道19号 博悦宴会艺术中心"
    $ws.Cells.Item(12, 5).Value = "2024.02.23 09:00-02.23 17:00"
    $ws.Cells.Item(12, 6).Value = 722
    $ws.Cells.Item(12, 7).Value = 65
    $ws.Cells.Item(12, 8).Value = "https://show.bilibili.com/platform/detail.html?id=80240"
    $ws.Cells.Item(12, 9).Value = "//i0.hdslb.com/bfs/openplatform/202312/Qwh83wl31703836740097.jpeg"

    $ws.Cells.Item(13, 2).Value = "2024.02.23"
    $ws.Cells.Item(13, 3).Value = "南昌·国乙only·突破次元计划（取消）"
    $ws.Cells.Item(13, 4).Value = "高处见美好生活公园 百家喜宴高新店"
    $ws.Cells.Item(13, 5).Value = "2024.02.23 10:00-02.23 21:00"
    $ws.Cells.Item(13, 6).Value = 304
    $ws.Cells.Item(13, 7).Value = "不可售"
    $ws.Cells.Item(13, 8).Value = "https://show.bilibili.com/platform/detail.html?id=80413"
    $ws.Cells.Item(13, 9).Value = "//i0.hdslb.com/bfs/openplatform/202401/XvmB77wb1704252353395.jpeg"

    $ws.Cells.Item(14, 2).Value = "2024.02.24"
    $ws.Cells.Item(14, 3).Value = "南昌·Cookie动漫嘉年华-赵路专场票"
    $ws.Cells.Item(14, 4).Value = "九龙大道1177号 南昌绿地国际博览中心"
    $ws.Cells.Item(14, 5).Value = "2024.02.24 11:00-02.24 17:00"
    $ws.Cells.Item(14, 6).Value = 352
    $ws.Cells.Item(14, 7).Value = "已售罄"
    $ws.Cells.Item(14, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81769"
    $ws.Cells.Item(14, 9).Value = "//i2.hdslb.com/bfs/openplatform/202402/DhCi2kWe1707123386859.png"

    $ws.Cells.Item(15, 2).Value = "2024.02.24"
    $ws.Cells.Item(15, 3).Value = "南昌·第一届Cookie动漫嘉年华"
    $ws.Cells.Item(15, 4).Value = "九龙大道1177号 南昌绿地国际博览中心"
    $ws.Cells.Item(15, 5).Value = "2024.02.24 09:00-02.24 17:00"
    $ws.Cells.Item(15, 6).Value = 4388
    $ws.Cells.Item(15, 7).Value = 65
    $ws.Cells.Item(15, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81033"
    $ws.Cells.Item(15, 9).Value = "//i1.hdslb.com/bfs/openplatform/202401/P994oBkz1705562167665.png"

    $ws.Cells.Item(16, 2).Value = "2024.02.24"
    $ws.Cells.Item(16, 3).Value = "宜春·融荟城难忘今宵汉文化节"
    $ws.Cells.Item(16, 4).Value = "宜阳大道239号 宜春融荟城"
    $ws.Cells.Item(16, 5).Value = "2024.02.24 14:00-02.24 18:00"
    $ws.Cells.Item(16, 6).Value = 19
    $ws.Cells.Item(16, 7).Value = 10
    $ws.Cells.Item(16, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81690"
    $ws.Cells.Item(16, 9).Value = "//i0.hdslb.com/bfs/openplatform/202402/ldtkc9Sp1706865634128.jpeg"

    $ws.Cells.Item(17, 2).Value = "2024.02.24"
    $ws.Cells.Item(17, 3).Value = "景德镇·陶溪川×次元文化元宵游园会（ 免费活动）"
    $ws.Cells.Item(17, 4).Value = "新厂西路315号 陶溪川发布大厅"
    $ws.Cells.Item(17, 5).Value = "2024.02.24 10:00-02.25 18:00"
    $ws.Cells.Item(17, 6).Value = 322
    $ws.Cells.Item(17, 7).Value = 30
    $ws.Cells.Item(17, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81207"
    $ws.Cells.Item(17, 9).Value = "//i1.hdslb.com/bfs/openplatform/202402/nIs2jtUn1707298876430.png"

    $ws.Cells.Item(18, 2).Value = "2024.03.02"
    $ws.Cells.Item(18, 3).Value = "南昌·meeting动漫游戏嘉年华"
    $ws.Cells.Item(18, 4).Value = "南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆"
    $ws.Cells.Item(18, 5).Value = "2024.03.02 09:00-03.03 17:00"
    $ws.Cells.Item(18, 6).Value = 1184
    $ws.Cells.Item(18, 7).Value = 60
    $ws.Cells.Item(18, 8).Value = "https://show.bilibili.com/platform/detail.html?id=79555"
    $ws.Cells.Item(18, 9).Value = "//i0.hdslb.com/bfs/openplatform/202402/l6GUtggC1706843695971.jpeg"

    $ws.Cells.Item(19, 2).Value = "2024.03.09"
    $ws.Cells.Item(19, 3).Value = "景德镇·江报国风动漫展 "
    $ws.Cells.Item(19, 4).Value = "迎宾大道与寺山路交叉口东200米 陶博城"
    $ws.Cells.Item(19, 5).Value = "2024.03.09 09:00-03.10 17:00"
    $ws.Cells.Item(19, 6).Value = 509
    $ws.Cells.Item(19, 7).Value = 45
    $ws.Cells.Item(19, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81362"
    $ws.Cells.Item(19, 9).Value = "//i0.hdslb.com/bfs/openplatform/202401/ae5G3ouV1706092057911.jpeg"

    $ws.Cells.Item(20, 2).Value = "2024.03.16"
    $ws.Cells.Item(20, 3).Value = "景德镇·原神X崩铁X崩坏动漫展only"
    $ws.Cells.Item(20, 4).Value = "陶阳南路188号 晨枫臻品酒店"
    $ws.Cells.Item(20, 5).Value = "2024.03.16 10:00-03.16 17:00"
    $ws.Cells.Item(20, 6).Value = 49
    $ws.Cells.Item(20, 7).Value = 55
    $ws.Cells.Item(20, 8).Value = "https://show.bilibili.com/platform/detail.html?id=80920"
    $ws.Cells.Item(20, 9).Value = "//i0.hdslb.com/bfs/openplatform/202401/IugBckTp1705469476482.png"

    $ws.Cells.Item(21, 2).Value = "2024.03.16"
    $ws.Cells.Item(21, 3).Value = "江西·ShiningStaR动漫游戏文化节5th"
    $ws.Cells.Item(21, 4).Value = "高新开发区紫阳大道666号 江西奥林匹克体育中心综合训练馆"
    $ws.Cells.Item(21, 5).Value = "2024.03.16 09:30-03.17 17:00"
    $ws.Cells.Item(21, 6).Value = 746
    $ws.Cells.Item(21, 7).Value = 60
    $ws.Cells.Item(21, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81792"
    $ws.Cells.Item(21, 9).Value = "//i2.hdslb.com/bfs/openplatform/202402/2l16aHBJ1707209383729.jpeg"

    $ws.Cells.Item(22, 2).Value = "2024.03.23"
    $ws.Cells.Item(22, 3).Value = "上饶·原×铁×崩only"
    $ws.Cells.Item(22, 4).Value = "五三东大道42号 回禾酒店"
    $ws.Cells.Item(22, 5).Value = "2024.03.23 10:00-03.23 17:00"
    $ws.Cells.Item(22, 6).Value = 26
    $ws.Cells.Item(22, 7).Value = 60
    $ws.Cells.Item(22, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81103"
    $ws.Cells.Item(22, 9).Value = "//i2.hdslb.com/bfs/openplatform/202401/pp6c5TsC1705647180602.jpeg"

    $ws.Cells.Item(23, 2).Value = "2024.03.23"
    $ws.Cells.Item(23, 3).Value = "南昌·AP动漫游戏嘉年华"
    $ws.Cells.Item(23, 4).Value = "八一桥街道青山南路118号 蓝海会展中心"
    $ws.Cells.Item(23, 5).Value = "2024.03.23 09:00-03.24 17:00"
    $ws.Cells.Item(23, 6).Value = 390
    $ws.Cells.Item(23, 7).Value = 60
    $ws.Cells.Item(23, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81232"
    $ws.Cells.Item(23, 9).Value = "//i2.hdslb.com/bfs/openplatform/202401/NZv97SmS1705912230957.jpeg"

    $ws.Cells.Item(24, 2).Value = "2024.03.23"
    $ws.Cells.Item(24, 3).Value = "南昌·原X穹X崩only"
    $ws.Cells.Item(24, 4).Value = "丰和北大道299号 新吉花园酒店"
    $ws.Cells.Item(24, 5).Value = "2024.03.23 10:00-03.23 17:00"
    $ws.Cells.Item(24, 6).Value = 52
    $ws.Cells.Item(24, 7).Value = 65
    $ws.Cells.Item(24, 8).Value = "https://show.bilibili.com/platform/detail.html?id=80807"
    $ws.Cells.Item(24, 9).Value = "//i0.hdslb.com/bfs/openplatform/202401/rY4v2Opx1705051458246.jpeg"

    $ws.Cells.Item(25, 2).Value = "2024.03.30"
    $ws.Cells.Item(25, 3).Value = "南昌·CM01动漫游戏博览会"
    $ws.Cells.Item(25, 4).Value = "怀玉山大道1315号 南昌绿地国际博览中心"
    $ws.Cells.Item(25, 5).Value = "2024.03.30 10:00-03.31 17:00"
    $ws.Cells.Item(25, 6).Value = 196
    $ws.Cells.Item(25, 7).Value = 55
    $ws.Cells.Item(25, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81691"
    $ws.Cells.Item(25, 9).Value = "//i2.hdslb.com/bfs/openplatform/202402/IYLaH7AS1706866218597.png"

    $ws.Cells.Item(26, 2).Value = "2024.03.30"
    $ws.Cells.Item(26, 3).Value = "鹰潭·原×铁×崩only"
    $ws.Cells.Item(26, 4).Value = "南站路24号 回禾酒店(鹰潭火车站店)"
    $ws.Cells.Item(26, 5).Value = "2024.03.30 10:00-03.30 17:00"
    $ws.Cells.Item(26, 6).Value = 16
    $ws.Cells.Item(26, 7).Value = 60
    $ws.Cells.Item(26, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81097"
    $ws.Cells.Item(26, 9).Value = "//i2.hdslb.com/bfs/openplatform/202401/q0AZaXAk1705646244207.jpeg"

    # Remove the now-obsolete trailing rows (old rows 27-31) and shrink
    # the sheets used range down to A1:I26.
    $ws.Range("A27:I31").Delete()
}

